# Add new columns I (I0) and J (IF) to the sheet, matching the style of the
# existing header cells and filling in the data for rows 2-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing headers (B1:H1) by copying
# the format from H1 (reuses the same style index rather than minting a
# new one).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-43) for columns I and J ---
$data = @{
    2  = @(9, 9)
    3  = @(7, 7)
    4  = @(9, 9)
    5  = @(5, 6)
    6  = @(6, 6)
    7  = @(4, 5)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(7, 7)
    11 = @(5, 6)
    12 = @(4, 5)
    13 = @(10, 10)
    14 = @(9, 9)
    15 = @(7, 8)
    16 = @(9, 9)
    17 = @(9, 9)
    18 = @(3, 4)
    19 = @(6, 7)
    20 = @(8, 9)
    21 = @(9, 9)
    22 = @(9, 9)
    23 = @(5, 6)
    24 = @(9, 9)
    25 = @(8, 8)
    26 = @(6, 6)
    27 = @(6, 6)
    28 = @(8, 8)
    29 = @(9, 9)
    30 = @(9, 9)
    31 = @(8, 8)
    32 = @(6, 6)
    33 = @(4, 5)
    34 = @(10, 10)
    35 = @(6, 6)
    36 = @(9, 9)
    37 = @(5, 5)
    38 = @(8, 8)
    39 = @(8, 9)
    40 = @(6, 6)
    41 = @(7, 7)
    42 = @(7, 7)
    43 = @(2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
